# Updates Fgf15-Fgfr4 LR-pairs sheet with refreshed TPM-derived values and two
# additional sending/target cluster combinations (rows 10-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @('Inflammatory-Mac', 'Fgf15', 'Fgfr4', 'ECs', 1, 0.3333333333333333, 0.018508, 0.055524, 0.2347866901774728, 0.3151800006811757, 3, 1, 0.363908, 1.091724, 0.01118972054940699, 0.01663265480083782, 0.006735209264, 0.06061688337599999, 0.00262719745180612, 0.005242280151457824),
    @('Inflammatory-Mac', 'Fgf15', 'Fgfr4', 'FAPs', 1, 0.3333333333333333, 0.018508, 0.055524, 0.2347866901774728, 0.3151800006811757, 1, 0.3333333333333333, 0.153566, 0.460698, 0.004721964413781051, 0.007018835164781924, 0.002842199528, 0.025579795752, 0.001108654395847463, 0.002212196472017027),
    @('Inflammatory-Mac', 'Fgf15', 'Fgfr4', 'Inflammatory-Mac', 1, 0.3333333333333333, 0.018508, 0.055524, 0.2347866901774728, 0.3151800006811757, 1, 0.3333333333333333, 0.06252866666666666, 0.187586, 0.001922679101110775, 0.002857913889838424, 0.001157280562666667, 0.010415525064, 0.0004514194624231975, 0.000900757301746016),
    @('Inflammatory-Mac', 'Fgf15', 'Fgfr4', 'MuSCs', 1, 0.3333333333333333, 0.018508, 0.055524, 0.2347866901774728, 0.3151800006811757, 2, 1, 31.927516, 63.85503199999999, 0.9817315966582778, 0.9728454303033116, 0.5909144661279999, 3.545486796767999, 0.2304975122220427, 0.3066214233856764),
    @('Inflammatory-Mac', 'Fgf15', 'Fgfr4', 'Resolving-Mac', 1, 0.3333333333333333, 0.018508, 0.055524, 0.2347866901774728, 0.3151800006811757, 1, 0.3333333333333333, 0.01411566666666667, 0.042347, 0.0004340392774233579, 0.0006451658412300904, 0.0002612527586666667, 0.002351274828, 0.0001019066453532521, 0.0002033433702783712),
    @('Resolving-Mac', 'Fgf15', 'Fgfr4', 'ECs', 1, 0.5, 0.060321, 0.120642, 0.7652133098225272, 0.6848199993188243, 3, 1, 0.363908, 1.091724, 0.01118972054940699, 0.01663265480083782, 0.021951294468, 0.131707766808, 0.008562523097600873, 0.01139037464938),
    @('Resolving-Mac', 'Fgf15', 'Fgfr4', 'FAPs', 1, 0.5, 0.060321, 0.120642, 0.7652133098225272, 0.6848199993188243, 1, 0.3333333333333333, 0.153566, 0.460698, 0.004721964413781051, 0.007018835164781924, 0.009263254686, 0.055579528116, 0.003613310017933587, 0.004806638692764898),
    @('Resolving-Mac', 'Fgf15', 'Fgfr4', 'Inflammatory-Mac', 1, 0.5, 0.060321, 0.120642, 0.7652133098225272, 0.6848199993188243, 1, 0.3333333333333333, 0.06252866666666666, 0.187586, 0.001922679101110775, 0.002857913889838424, 0.003771791702, 0.022630750212, 0.001471259638687578, 0.001957156588092408),
    @('MuSCs', 'Fgf15', 'Fgfr4', 'MuSCs', 1, 0.5, 0.060321, 0.120642, 0.7652133098225272, 0.6848199993188243, 2, 1, 31.927516, 63.85503199999999, 0.9817315966582778, 0.9728454303033116, 1.925899692636, 7.703598770544, 0.7512340844362351, 0.6662240069176352),
    @('MuSCs', 'Fgf15', 'Fgfr4', 'Resolving-Mac', 1, 0.5, 0.060321, 0.120642, 0.7652133098225272, 0.6848199993188243, 1, 0.3333333333333333, 0.01411566666666667, 0.042347, 0.0004340392774233579, 0.0006451658412300904, 0.000851471129, 0.005108826774000001, 0.0003321326320701058, 0.0004418224709517192)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $excelRow).Value = $rowVals[$c]
    }
}